$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics")
$ws.Range("B2").Value = 406461.36
$ws.Range("D20").Value = "test"
